$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) stays unchanged.
# Data rows: row 2 is updated in place, rows 3-7 are newly added.
$rows = @{
    2 = @("ECs",  "Ntf3", "Ntrk1", "FAPs", 2, 0.6666666666666666, 7.310333666666668,  21.931001,           0.357182641950452,  0.357182641950452,  3, 1,                  0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 0.8856639811618892, 7.970975830457001,  0.2330525622668684,  0.2330525622668684)
    3 = @("ECs",  "Ntf3", "Ntrk1", "sCs",  2, 0.6666666666666666, 7.310333666666668,  21.931001,           0.357182641950452,  0.357182641950452,  1, 0.3333333333333333, 0.064529,            0.193587, 0.3475255096545335, 0.3475255096545335, 0.4717285211763334, 4.245556690587001,  0.1241300796835836,  0.1241300796835836)
    4 = @("FAPs", "Ntf3", "Ntrk1", "FAPs", 3, 1,                  10.17625966666667,  30.528779,           0.4972116839874968, 0.4972116839874969, 3, 1,                  0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 1.232877603222556,  11.095898429003,    0.3244179401035531,  0.3244179401035531)
    5 = @("FAPs", "Ntf3", "Ntrk1", "sCs",  3, 1,                  10.17625966666667,  30.528779,           0.4972116839874968, 0.4972116839874969, 1, 0.3333333333333333, 0.064529,            0.193587, 0.3475255096545335, 0.3475255096545335, 0.6566638600303334, 5.909974740273,     0.1727937438839437,  0.1727937438839437)
    6 = @("sCs",  "Ntf3", "Ntrk1", "FAPs", 2, 0.6666666666666666, 2.980061000000001,  8.940183000000001,  0.1456056740620512, 0.1456056740620512, 3, 1,                  0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 0.3610413436256668, 3.249372092631001,  0.09500398797504492, 0.09500398797504492)
    7 = @("sCs",  "Ntf3", "Ntrk1", "sCs",  2, 0.6666666666666666, 2.980061000000001,  8.940183000000001,  0.1456056740620512, 0.1456056740620512, 1, 0.3333333333333333, 0.064529,            0.193587, 0.3475255096545335, 0.3475255096545335, 0.1923003562690001, 1.730703206421,     0.05060168608700623, 0.05060168608700623)
}

foreach ($r in 2..7) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
